# Matriz.xlsx edit script
# Inserts 5 new test rows (Clima, Deportes, Espectaculos, Tecnologia,
# Informacion sobre covid-19, and a new "Validar Hipervinculo del logo" row)
# before the existing Login/Noticias rows, which shift down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Push the existing rows 4-10 down to rows 9-15, leaving 5 blank rows
#    (4-8) to be populated below. Row insert preserves each shifted row's
#    own height/content/styles automatically.
$ws.Rows("4:8").Insert()

# 2. Update row 3 (was the generic "menu" case) and populate the 5 new
#    rows (4-8) with the per-section test cases + the new logo test case.

# Row 3 - PPDN-002 (Clima)
$ws.Range("A3").Value = "PPDN-002"
$ws.Range("B3").Value = "Pagina principal"
$ws.Range("C3").Value = "Validar Hipervinculos`ndel menu"
$ws.Range("D3").Value = "Acceso a la pagina`nprincipal"
$ws.Range("E3").Value = "1.Entrar a la pagina prinicpal.`n2.Presionar la seccion `"Clima`""
$ws.Range("F3").Value = "1.Mandar a la seccion de Clima"
$ws.Range("G3").Value = "Esperado"
$ws.Range("H3").Value = "Correcto"

# Row 4 - PPDN-003 (Deportes)
$ws.Range("A4").Value = "PPDN-003"
$ws.Range("B4").Value = "Pagina principal"
$ws.Range("C4").Value = "Validar Hipervinculos`ndel menu"
$ws.Range("D4").Value = "Acceso a la pagina`nprincipal"
$ws.Range("E4").Value = "1.Entrar a la pagina prinicpal.`n2.Presionar la seccion `"Deportes`""
$ws.Range("F4").Value = "1.Mandar a la seccion de Deportes"
$ws.Range("G4").Value = "Esperado"
$ws.Range("H4").Value = "Correcto"
$ws.Rows(4).RowHeight = 60

# Row 5 - PPDN-004 (Espectaculos)
$ws.Range("A5").Value = "PPDN-004"
$ws.Range("B5").Value = "Pagina principal"
$ws.Range("C5").Value = "Validar Hipervinculos`ndel menu"
$ws.Range("D5").Value = "Acceso a la pagina`nprincipal"
$ws.Range("E5").Value = "1.Entrar a la pagina prinicpal.`n2.Presionar la seccion `"Espectaculos`""
$ws.Range("F5").Value = "1.Mandar a la seccion de Espectaculos"
$ws.Range("G5").Value = "Esperado"
$ws.Range("H5").Value = "Correcto"

# Row 6 - PPDN-005 (Tecnologia)
$ws.Range("A6").Value = "PPDN-005"
$ws.Range("B6").Value = "Pagina principal"
$ws.Range("C6").Value = "Validar Hipervinculos`ndel menu"
$ws.Range("D6").Value = "Acceso a la pagina`nprincipal"
$ws.Range("E6").Value = "1.Entrar a la pagina prinicpal.`n2.Presionar la seccion `"Tecnologia`""
$ws.Range("F6").Value = "1.Mandar a la seccion de Tecnologia"
$ws.Range("G6").Value = "Esperado"
$ws.Range("H6").Value = "Correcto"
$ws.Rows(6).RowHeight = 60

# Row 7 - PPDN-006 (Informacion sobre covid-19)
$ws.Range("A7").Value = "PPDN-006"
$ws.Range("B7").Value = "Pagina principal"
$ws.Range("C7").Value = "Validar Hipervinculos`ndel menu"
$ws.Range("D7").Value = "Acceso a la pagina`nprincipal"
$ws.Range("E7").Value = "1.Entrar a la pagina prinicpal.`n2.Presionar la seccion `"Informacion sobre covid-19`""
$ws.Range("F7").Value = "1.Mandar a la seccion de Informacion sobre covid-19"
$ws.Range("G7").Value = "Esperado"
$ws.Range("H7").Value = "Correcto"

# Row 8 - PPDN-007 (Validar Hipervinculo del logo) - brand-new test case
$ws.Range("A8").Value = "PPDN-007"
$ws.Range("B8").Value = "Pagina principal"
$ws.Range("C8").Value = "Validar Hipervinculo del logo"
$ws.Range("D8").Value = "Acceso a la pagina`nprincipal"
$ws.Range("E8").Value = "1.Entrar a la pagina prinicpal.`n2.Presionar el logotipo "
$ws.Range("F8").Value = "1.Regresar a la pagina principal"
$ws.Range("G8").Value = "Esperado"
$ws.Range("H8").Value = "Correcto"

# 3. Fix up row heights on the rows that were shifted down but whose
#    heights differ from what they had before the insert.
$ws.Rows(10).RowHeight = 75
$ws.Rows(12).RowHeight = 90

# 4. Update the view: scroll near the top and select I11 (matches the
#    authored selection after the edit).
$ws.Range("A3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I11").Select()
